# update dq_msg and umlauts

$wb = $excel.ActiveWorkbook

# --- Sheet "DQ_Report" (sheet1) ---
$wsDQ = $wb.Worksheets.Item("DQ_Report")

# Header rename: ICD_Primärkode -> ICD_primaerkode
$wsDQ.Range("B1").Value = "ICD_primaerkode"

# Remove row 16 (P_20085770 / J09) entirely, shrinking used range to A1:C15
$wsDQ.Rows.Item(16).Delete()

# --- Sheet "Statistik" (sheet2) ---
$wsStat = $wb.Worksheets.Item("Statistik")

# Header renames
$wsStat.Range("F1").Value = "K2_icdRd_no"
$wsStat.Range("G1").Value = "K3_rd_no"

# Value updates
$wsStat.Range("E2").Value = 97.7
$wsStat.Range("G2").Value = 297

# --- Sheet "Projectathon" (sheet3) ---
$wsProj = $wb.Worksheets.Item("Projectathon")

# Header rename: ICD_Primärkode -> ICD_primaerkode
$wsProj.Range("D1").Value = "ICD_primaerkode"
